$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.55

# Row 3
$ws.Range("G3").Value = 2.5

# Row 6
$ws.Range("N6").Value = 1.9
$ws.Range("O6").Value = 1.95

# Row 9
$ws.Range("G9").Value = 2.8
$ws.Range("I9").Value = 3.25
$ws.Range("K9").Value = 4.3
$ws.Range("L9").Value = 1.65
$ws.Range("M9").Value = 2.1
$ws.Range("N9").Value = 2.95
$ws.Range("O9").Value = 1.35
$ws.Range("T9").Value = 5.8
$ws.Range("U9").Value = 12.5
$ws.Range("W9").Value = 37
$ws.Range("X9").Value = 32
$ws.Range("Z9").Value = 4.3
$ws.Range("AE9").Value = 6.6
$ws.Range("AF9").Value = 15.5

# Row 10
$ws.Range("G10").Value = 2.15
$ws.Range("J10").Value = 1.07
$ws.Range("L10").Value = 1.4
$ws.Range("AD10").Value = 501

# Row 13
$ws.Range("L13").Value = 1.3
$ws.Range("M13").Value = 3.25
$ws.Range("R13").Value = 1.99
$ws.Range("S13").Value = 1.73

# Row 15
$ws.Range("J15").Value = 1.11
$ws.Range("K15").Value = 6.5
$ws.Range("P15").Value = 1.53
$ws.Range("Q15").Value = 2.38
$ws.Range("U15").Value = 12
$ws.Range("V15").Value = 11

# Row 16
$ws.Range("J16").Value = 1.1
$ws.Range("K16").Value = 7

# Row 17
$ws.Range("K17").Value = 8.5
$ws.Range("Y17").Value = 29
$ws.Range("Z17").Value = 8.5
$ws.Range("AD17").Value = 451

# Row 18
$ws.Range("J18").Value = 1.1
$ws.Range("K18").Value = 7
$ws.Range("L18").Value = 1.5
$ws.Range("M18").Value = 2.5
$ws.Range("N18").Value = 2.5
$ws.Range("O18").Value = 1.5
$ws.Range("P18").Value = 1.53
$ws.Range("Q18").Value = 2.38
$ws.Range("R18").Value = 2.05
$ws.Range("S18").Value = 1.7
$ws.Range("T18").Value = 7
$ws.Range("Z18").Value = 7
$ws.Range("AC18").Value = 67
$ws.Range("AD18").Value = 501
$ws.Range("AI18").Value = 26

# Row 19
$ws.Range("J19").Value = 1.06
$ws.Range("K19").Value = 10
$ws.Range("N19").Value = 2.08
$ws.Range("O19").Value = 1.73

# Row 22
$ws.Range("J22").Value = 1.13
$ws.Range("K22").Value = 6
$ws.Range("R22").Value = 2.2
$ws.Range("S22").Value = 1.62
$ws.Range("T22").Value = 6.5
$ws.Range("Z22").Value = 6
$ws.Range("AC22").Value = 81
$ws.Range("AG22").Value = 12
$ws.Range("AI22").Value = 29

# Row 23
$ws.Range("G23").Value = 2.7
$ws.Range("H23").Value = 3.6
$ws.Range("I23").Value = 2.35
$ws.Range("N23").Value = 1.67
$ws.Range("O23").Value = 2.15
$ws.Range("T23").Value = 11
$ws.Range("U23").Value = 15
$ws.Range("V23").Value = 10
$ws.Range("AA23").Value = 7
$ws.Range("AH23").Value = 23
$ws.Range("AJ23").Value = 23

# Row 25
$ws.Range("G25").Value = 3.1
$ws.Range("H25").Value = 3
$ws.Range("I25").Value = 2.38
$ws.Range("N25").Value = 2.5
$ws.Range("O25").Value = 1.5
$ws.Range("P25").Value = 1.53
$ws.Range("Q25").Value = 2.38
$ws.Range("R25").Value = 2.1
$ws.Range("S25").Value = 1.67
$ws.Range("T25").Value = 7.5
$ws.Range("Z25").Value = 7
$ws.Range("AH25").Value = 23
$ws.Range("AI25").Value = 23
$ws.Range("AJ25").Value = 41

# Row 26
$ws.Range("L26").Value = 1.22
$ws.Range("M26").Value = 4

# Row 27
$ws.Range("P27").Value = 1.4
$ws.Range("Q27").Value = 2.75
$ws.Range("R27").Value = 1.83
$ws.Range("S27").Value = 1.83
$ws.Range("T27").Value = 8.5
$ws.Range("Z27").Value = 10
$ws.Range("AA27").Value = 7
$ws.Range("AD27").Value = 301
$ws.Range("AE27").Value = 7.5
$ws.Range("AG27").Value = 9.5
$ws.Range("AJ27").Value = 29

# Row 29
$ws.Range("G29").Value = 2.55
$ws.Range("I29").Value = 2.63
$ws.Range("N29").Value = 2.03
$ws.Range("O29").Value = 1.83
$ws.Range("U29").Value = 13
$ws.Range("W29").Value = 26
$ws.Range("AB29").Value = 13
$ws.Range("AD29").Value = 201
$ws.Range("AE29").Value = 8.5

# Row 30
$ws.Range("G30").Value = 2.7
$ws.Range("I30").Value = 2.55
$ws.Range("K30").Value = 10
$ws.Range("X30").Value = 21
$ws.Range("Y30").Value = 29
$ws.Range("AH30").Value = 26

# Row 31
$ws.Range("K31").Value = 9.5
$ws.Range("T31").Value = 7
$ws.Range("AF31").Value = 17

# Row 34
$ws.Range("G34").Value = 3.1
$ws.Range("I34").Value = 2.32
$ws.Range("K34").Value = 5.9
$ws.Range("L34").Value = 1.52
$ws.Range("M34").Value = 2.42
$ws.Range("N34").Value = 2.5
$ws.Range("P34").Value = 1.57
$ws.Range("Q34").Value = 2.32
$ws.Range("R34").Value = 2.15
$ws.Range("S34").Value = 1.62
$ws.Range("T34").Value = 7.1
$ws.Range("U34").Value = 15
$ws.Range("X34").Value = 37
$ws.Range("Z34").Value = 5.9
$ws.Range("AE34").Value = 6
$ws.Range("AF34").Value = 10.25
$ws.Range("AG34").Value = 10.75
$ws.Range("AH34").Value = 25
$ws.Range("AI34").Value = 26

# Row 38
$ws.Range("J38").Value = 1.01
$ws.Range("K38").Value = 34
$ws.Range("L38").Value = 1.07

# Row 39
$ws.Range("J39").Value = 1.03
$ws.Range("L39").Value = 1.17

# Row 40
$ws.Range("G40").Value = 1.85
$ws.Range("H40").Value = 3.5
$ws.Range("I40").Value = 3.9
$ws.Range("J40").Value = 1.05
$ws.Range("M40").Value = 3.75
$ws.Range("N40").Value = 1.88
$ws.Range("O40").Value = 1.98
$ws.Range("R40").Value = 1.75
$ws.Range("S40").Value = 2
$ws.Range("T40").Value = 8
$ws.Range("Y40").Value = 23
$ws.Range("Z40").Value = 11
$ws.Range("AA40").Value = 7
$ws.Range("AC40").Value = 41
$ws.Range("AD40").Value = 201
$ws.Range("AE40").Value = 12
$ws.Range("AI40").Value = 29
$ws.Range("AJ40").Value = 34
